# Demo AutomationFrame work release 1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AssetProfile")
$ws.Activate()

# Renumber the Demo-3xx asset ids in column A (rows 2-7)
$ws.Range("A2").Value = "Demo-322"
$ws.Range("A3").Value = "Demo-323"
$ws.Range("A4").Value = "Demo-324"
$ws.Range("A5").Value = "Demo-325"
$ws.Range("A6").Value = "Demo-326"
$ws.Range("A7").Value = "Demo-327"

# Update the active selection to the asset id column
$ws.Range("A2:A7").Select()
